# This script applies the updated 'Ventas Objetivo' related recalculation
# (columns L 'Diferencia Stock', R 'uds. Objetivo semana pasada',
# T 'Tendencia Consumo' = MAX(S - R, 0)) plus the C57 total adjustment,
# as described in the commit: 'Hemos cambiado la formula de Ventas objetivo'.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

# Row 5
$ws.Range("L5").Value = 0
$ws.Range("R5").Value = 5
$ws.Range("T5").Value = 0

# Row 6
$ws.Range("L6").Value = 0
$ws.Range("R6").Value = 2

# Row 7
$ws.Range("R7").Value = 6

# Row 8
$ws.Range("L8").Value = 0
$ws.Range("R8").Value = 7
$ws.Range("T8").Value = 0

# Row 9
$ws.Range("L9").Value = 0
$ws.Range("R9").Value = 14
$ws.Range("T9").Value = 0

# Row 10
$ws.Range("R10").Value = 1

# Row 11
$ws.Range("R11").Value = 1
$ws.Range("T11").Value = 1

# Row 12
$ws.Range("R12").Value = 2
$ws.Range("T12").Value = 1

# Row 13
$ws.Range("R13").Value = 6
$ws.Range("T13").Value = 0

# Row 16
$ws.Range("L16").Value = 0

# Row 20
$ws.Range("L20").Value = 0

# Row 23
$ws.Range("R23").Value = 3

# Row 25
$ws.Range("R25").Value = 2
$ws.Range("T25").Value = 0

# Row 28
$ws.Range("R28").Value = 5
$ws.Range("T28").Value = 0

# Row 31
$ws.Range("L31").Value = 0

# Row 32
$ws.Range("R32").Value = 1
$ws.Range("T32").Value = 7

# Row 34
$ws.Range("L34").Value = 0

# Row 35
$ws.Range("L35").Value = 0
$ws.Range("R35").Value = 13
$ws.Range("T35").Value = 0

# Row 36
$ws.Range("R36").Value = 7

# Row 37
$ws.Range("L37").Value = 0
$ws.Range("R37").Value = 6
$ws.Range("T37").Value = 0

# Row 38
$ws.Range("R38").Value = 2

# Row 39
$ws.Range("R39").Value = 3

# Row 40
$ws.Range("L40").Value = 0
$ws.Range("R40").Value = 6
$ws.Range("T40").Value = 10

# Row 41
$ws.Range("L41").Value = 0
$ws.Range("R41").Value = 4
$ws.Range("T41").Value = 0

# Row 42
$ws.Range("R42").Value = 1

# Row 57
$ws.Range("C57").Value = 0
